# The DMD sheet ("Demand") has a small timeslice-share table in B23:D33
# (headers "~FI_T:" / "Timeslice" / "COM_FR" plus 8 data rows referencing
# SEC_Comm!$C$8 and the 1S1W1D..2S2W2D timeslice codes) that is being
# wiped out, leaving the cell formatting in place but with empty content.
# D38 holds =SUM(D26:D33), which recalculates from 1 down to 0 once the
# D26:D33 inputs are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DMD")

# Clear all text/formulas/values from the timeslice table, keeping styles.
$ws.Range("B23:D33").ClearContents()

# The label row (B25:D25) loses its two-line wrapped text, so its manual
# row height collapses back to the single-line default for this sheet.
$ws.Rows.Item(25).RowHeight = 13.5

# Leave the selection where the now-cleared table used to be.
$ws.Range("B23:D33").Select()
